# GPLIM-3541: add Material Type as required header for Manifest uploads
#
# Adds a new "Material Type" column (G) to the manifest worksheet:
#   - G1 header "Material Type" styled bold/white-on-black, centered
#   - G2:G24 data cells "DNA:Genomic", centered (reuses the existing
#     centered style used by column C)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell (G1) ------------------------------------------------
# Start from the F1 header's style (bold font + colored fill, centered,
# no quote-prefix baggage) so the new style we derive only needs its
# font color and fill color changed.
$headerSrc = $ws.Range("F1")
$headerDst = $ws.Range("G1")
$headerSrc.Copy()
$headerDst.PasteSpecial(-4122)

$headerDst.Value2 = "Material Type"
$headerDst.Interior.Color = 0
$headerDst.Interior.PatternColor = 0
$headerDst.Font.Color = 16777215

# --- Data cells (G2:G24) ---------------------------------------------
for ($r = 2; $r -le 24; $r++) {
    $cell = $ws.Range("G$r")
    $cell.Value2 = "DNA:Genomic"
    $cell.HorizontalAlignment = -4108
}

# --- Selection matches the authored workbook -------------------------
$ws.Range("G1:G24").Select() | Out-Null
